$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Business Requirement Document (BRD) Error Boundary in React",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Test Business Requirement Document (BRD) Error Boundary in React",
    2
)
